# Insert a new weekly price observation for "Arándano (blue)" / Vega Central
# Mapocho de Santiago at row 174, pushing the existing rows 174-280 down to
# 175-281 (dimension grows from A1:T280 to A1:T281).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 174..280 down to 175..281, leaving a blank row 174 to fill in.
$ws.Rows(174).Insert()

# Populate the newly inserted row 174 with the new record.
$ws.Cells.Item(174, 1).Value  = 9
$ws.Cells.Item(174, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(174, 3).Value  = "Metropolitana"
$ws.Cells.Item(174, 4).Value  = 44960
$ws.Cells.Item(174, 5).Value  = 13
$ws.Cells.Item(174, 6).Value  = "Fruta"
$ws.Cells.Item(174, 7).Value  = 100101
$ws.Cells.Item(174, 8).Value  = "Berries"
$ws.Cells.Item(174, 9).Value  = 100101001
$ws.Cells.Item(174, 10).Value = "Arándano (blue)"
$ws.Cells.Item(174, 11).Value = "Sin especificar"
$ws.Cells.Item(174, 12).Value = "Primera"
$ws.Cells.Item(174, 13).Value = 350
$ws.Cells.Item(174, 14).Value = 3000
$ws.Cells.Item(174, 15).Value = 3000
$ws.Cells.Item(174, 16).Value = 3000
$ws.Cells.Item(174, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(174, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(174, 19).Value = 1500
$ws.Cells.Item(174, 20).Value = 2
